# Apply cryptos list update (prices/volumes) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.499.47"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.253.30"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.861"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "2.268.17"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "42.270.40"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000102"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +25.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.31%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0825"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.19%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +10.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  +6.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.69%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.202"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("E51").Value = "  +0.71%  "
